# LOM3091.xlsx edit script
# Reconstructs the target state described by the diff: the sheet's label
# column (A) shifts up by one logical slot starting at row 13 (the
# "Objectives:"/"Docentes responsáveis:" pair loses its paired long text and
# the remaining label/value pairs cascade upward), a brand-new "Semestral"
# value appears, several long paragraph values are removed, and the last
# row (25) is dropped after its content is folded into row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Named constants for PasteSpecial (Excel XlPasteType enum)
$xlPasteAll     = -4104
$xlPasteFormats = -4122
$xlPasteValues  = -4163

# ---------------------------------------------------------------------
# Row 10: Objetivos: description cell -> reuse the "Durval" credential
# text that currently lives in B13/C13 (it will itself be overwritten
# further down), dropping the long "Apresentar e discutir..." paragraph.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"

# ---------------------------------------------------------------------
# Row 13: gains a new A13 label "Programa resumido:" and its B/C values
# switch from the old "Durval" text to the brand new "Semestral" value.
# ---------------------------------------------------------------------
$ws.Range("A14").Copy()
$ws.Range("A13").PasteSpecial($xlPasteFormats)
$ws.Range("A13").Value = "Programa resumido:"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Rows(13).RowHeight = 60

# ---------------------------------------------------------------------
# Row 14: label becomes "Short syllabus:"; its long descriptive
# paragraph (B14/C14) is removed entirely.
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# ---------------------------------------------------------------------
# Row 15: label becomes "Programa:"; B15/C15 are brand-new cells that
# take the date text "01/01/2012" already used by row 8 (copied so the
# value stays plain text instead of being reinterpreted as a date).
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Programa:"

$ws.Range("B19").Copy()
$ws.Range("B15").PasteSpecial($xlPasteFormats)
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial($xlPasteValues)

$ws.Range("C19").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial($xlPasteValues)

$ws.Rows(15).RowHeight = 120

# ---------------------------------------------------------------------
# Row 16: label becomes "Syllabus:"; the long numbered programa list in
# B16/C16 is removed entirely.
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# ---------------------------------------------------------------------
# Row 17: label becomes "Avaliação:"; it no longer needs the tall
# 120pt row height since it holds no B/C text any more.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows(17).AutoFit()

# ---------------------------------------------------------------------
# Row 18: label becomes "Método:"; gains new B18/C18 cells carrying the
# "Durval" credential text again (duplicated, matching the diff).
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Método:"

$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("B18").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("C18").Value = "6495737 - Durval Rodrigues Junior"

$ws.Rows(18).RowHeight = 60

# ---------------------------------------------------------------------
# Row 19: label becomes "Critério:"; B19/C19 keep their existing
# "Aplicação de duas provas..." text unchanged.
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Critério:"

# ---------------------------------------------------------------------
# Row 20: label becomes "Norma de recuperação:"; B20/C20 keep their
# existing "A Nota final..." text unchanged.
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Norma de recuperação:"

# ---------------------------------------------------------------------
# Row 21: label becomes "Bibliografia:"; B21/C21 keep their existing
# "A recuperação será feita..." text unchanged, but the row grows from
# 60pt to 120pt.
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

# ---------------------------------------------------------------------
# Row 22: label becomes "Requisitos:"; the long bibliography list in
# B22/C22 is removed, and the row height reverts to the default.
# ---------------------------------------------------------------------
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows(22).AutoFit()

# ---------------------------------------------------------------------
# Row 23: loses its A23 "Requisitos:" label; gains new B23/C23 cells
# carrying the "LOB1053 ... (Requisito)" text that currently sits in
# row 24 (copied before row 24 itself is overwritten below).
# ---------------------------------------------------------------------
$ws.Range("B19").Copy()
$ws.Range("B23").PasteSpecial($xlPasteFormats)
$ws.Range("B24").Copy()
$ws.Range("B23").PasteSpecial($xlPasteValues)

$ws.Range("C19").Copy()
$ws.Range("C23").PasteSpecial($xlPasteFormats)
$ws.Range("C24").Copy()
$ws.Range("C23").PasteSpecial($xlPasteValues)

$ws.Range("A23").Clear()
$ws.Rows(23).RowHeight = 30

# ---------------------------------------------------------------------
# Row 24: B24/C24 switch from the "LOB1053..." text to the
# "LOM3013 ... (Requisito)" text currently in row 25.
# ---------------------------------------------------------------------
$ws.Range("B25").Copy()
$ws.Range("B24").PasteSpecial($xlPasteValues)
$ws.Range("C25").Copy()
$ws.Range("C24").PasteSpecial($xlPasteValues)

# ---------------------------------------------------------------------
# Row 25 no longer exists in the target sheet (its content was folded
# into row 24 above), so the whole row is removed.
# ---------------------------------------------------------------------
$ws.Rows(25).Delete()
